# 1. type check bug fix (base type check is unprecise)
#
# Sheet1: the "type_object" example value had key2 typed as an array
# ([100,200]) even though the header documents {key1:int, key2:int} -
# fix the example so key2 is a plain int.
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("H5").Value = '{"key1":10100,"key2":20100}'

# Add a new "TestCase" sheet right after Sheet1 demonstrating the
# imprecise-type-check edge case (an int column fed a near-integer
# float inside a json object).
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("B1").Value = "TestCase"
$ws2.Range("A2").Value = "id"
$ws2.Range("B2").Value = "type"
$ws2.Range("A3").Value = "*int"
$ws2.Range("B3").Value = "{key1:int, key2:int}"
$ws2.Range("A4").Value = 1
$ws2.Range("B4").Value = '{"key1":10100,"key2":1.0000000001}'

# Restore selections: Sheet1 cursor moves to H6, Sheet2 becomes the
# active (selected) sheet with its cursor on B5.
[void]$ws1.Range("H6").Select()
[void]$ws2.Range("B5").Select()
[void]$ws2.Activate()
